# Add FunctionKey "clear" test rows to the PMTestData sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PMTestData")

# row, TestCase (A), Data (B), TestFlag (C), row height
$rows = @(
    @{ Row = 35; A = 'test_clearTNS_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function TNS --display-text "TNS-10001" --key 1 --key-sequence 10001,10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 101.5 },
    @{ Row = 36; A = 'test_clearMNS_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function MNS --display-text "MNS-10001" --key 1 --monitored-dir 10001 --alert-type 0,10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 101.5 },
    @{ Row = 37; A = 'test_clearEDN_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function EDN --display-text "EDN-10001" --key 1 --line-dir 10001,10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 87 },
    @{ Row = 38; A = 'test_clearMOI_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function XML --display-text "MOI" --key 1 --xml-on-demand-uri http://$$PROXYURL$$:22222/StreamingMenu?user=$$SIPUSERNAME$$,10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 116 },
    @{ Row = 39; A = 'test_clearPGM_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function PGM --display-text "PGM" --key 1,10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 87 },
    @{ Row = 40; A = 'test_clearREC_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function REC --display-text "REC" --key 1 --record-on-demand-uri ''http://149.13.0.80:80//nrj.ogg'',10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 101.5 },
    @{ Row = 41; A = 'test_clearDMN_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function DMN --display-text "DMN-10001" --key 1 --monitored-dir 10001 --alert-type 0,10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 101.5 },
    @{ Row = 42; A = 'test_clearGMA_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function GMA --display-text "GMA-10001" --key 1 --monitored-dir 10001,10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 101.5 },
    @{ Row = 43; A = 'test_clearMCT_funcKey'; B = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function MCT --display-text "MCT" --key 1,10000,10001,extension_key -e -d 10000 --key 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'; Height = 87 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 2).WrapText = $true
    $ws.Cells.Item($r.Row, 3).Value = "Y"
    $ws.Rows.Item($r.Row).RowHeight = $r.Height
}

# The first new row keeps an (empty) wrapped D cell like the rows above it.
$ws.Cells.Item(35, 4).WrapText = $true

# Update the view so the newly added rows are visible / selected.
$ws.Activate()
[void]$ws.Range("D41").Select()
try {
    $excel.ActiveWindow.ScrollRow = 41
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
